$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Bentleigh'
$ws.Range("B2").Value = 'Coles - Bentleigh  5/7 Vickery St, Bentleigh'
$ws.Range("C2").Value = '24/12/20 12:30pm-1:00pm'
$ws.Range("B3").Value = 'Il Centro Deli  5/284/292 Centre Rd, Bentleigh VIC 3204'
$ws.Range("C3").Value = '22/12/20 12:00pm-12:30pm'
$ws.Range("D3").Value = 'Case shopped in store'
$ws.Range("A4").Value = 'Black Rock'
$ws.Range("B4").Value = 'Woolworths Metro  40 Bluff Road, Black Rock VIC 3193'
$ws.Range("C4").Value = '30/12/20 5:30pm-5:55pm'
$ws.Range("D4").Value = 'Case shopped'
$ws.Range("A5").Value = 'Box Hill South'
$ws.Range("B5").Value = 'Bunnings  259 Middleborough Road, Box Hill South VIC 3128'
$ws.Range("C5").Value = '30/12/20 12:00pm-12:40pm'
$ws.Range("A6").Value = 'Brighton'
$ws.Range("B6").Value = 'Brighton Beach  Brighton, VIC 3186'
$ws.Range("C6").Value = '26/12/20 12:00pm-3:00pm'
$ws.Range("D6").Value = 'Case attended beach'
$ws.Range("C7").Value = '29/12/20 12:00pm-3:00pm'
$ws.Range("A8").Value = 'Burwood East'
$ws.Range("B8").Value = 'Coles Burwood, Burwood Highway & Blackburn Road'
$ws.Range("C8").Value = '28/12/20 6.30pm - 7pm'
$ws.Range("D8").Value = 'Case shopped in store'
$ws.Range("B9").Value = 'Kmart Burwood, 172 Burwood Highway'
$ws.Range("C9").Value = '28/12/20 6.15pm - 6.30pm'
$ws.Range("A10").Value = 'Camberwell'
$ws.Range("B10").Value = 'Coles Middle Camberwell, 751 Riversdale Road'
$ws.Range("C10").Value = '28/12/20 12pm - 12.30pm'
$ws.Range("B11").Value = 'Fu Lin Asian Grocery Supermarket  1397 Toorak Road, Camberwell VIC 3124'
$ws.Range("C11").Value = '30/12/20 2:30pm-2:45pm'
$ws.Range("D11").Value = 'Case shopped'
$ws.Range("A12").Value = 'Cape Schank'
$ws.Range("B12").Value = 'National Golf Club  The Cups Drive, Cape Schanck VIC 3939'
$ws.Range("C12").Value = '30/12/20 11.40am-1.40pm'
$ws.Range("D12").Value = 'Case attended course'
$ws.Range("A13").Value = 'Cheltenham'
$ws.Range("B13").Value = 'Aldi Cheltenham  280/282 Bay Road, Cheltenham VIC 3192'
$ws.Range("C13").Value = '29/12/20 2:00pm-2:30pm'
$ws.Range("D13").Value = 'Case shopped in store'
$ws.Range("B14").Value = 'Angus and Cootes Jeweller  Southland Shopping Centre, Shop 2096/1239, Nepean Hwy, Cheltenham VIC 3192'
$ws.Range("C14").Value = '28/12/2020 2:30pm-2:50pm'
$ws.Range("B15").Value = 'Bodero Southland Shopping Centre, 1239 Nepean Hwy'
$ws.Range("C15").Value = '22/12/20 6.45pm - 7pm'
$ws.Range("B16").Value = 'Chemist Warehouse Cheltenham, 326/330 Charman Rd'
$ws.Range("C16").Value = '30/12/20, 3.30pm - 3.45pm'
$ws.Range("B17").Value = 'Coles, Westfield Southland'
$ws.Range("C17").Value = '22/12/20 11:50am-12:10pm'
$ws.Range("B18").Value = 'Cotton On, Southland Shopping Centre 1239 Nepean Hwy'
$ws.Range("C18").Value = '22/12/20 12.15pm - 12.45pm'
$ws.Range("D18").Value = 'Case visited venue'
$ws.Range("B19").Value = 'Honey Birdette Southland  Shop 2209/1239, Southland Shopping Centre, Cheltenham VIC 3192'
$ws.Range("C19").Value = '22/12/2020 3:50pm-4:05pm'
$ws.Range("D19").Value = 'Case shopped in store'
$ws.Range("B20").Value = 'Kmart Southland Shopping Centre, 1239 Nepean Highway'
$ws.Range("C20").Value = '22/12/20 6.30pm - 6.45pm'
$ws.Range("C21").Value = '28/12/20 2.30pm-3pm'
$ws.Range("B22").Value = 'Mecca Southland  Shop 2011/2013, Southland Shopping Centre, Cheltenham VIC 3192'
$ws.Range("C22").Value = '22/12/2020 3:30pm-3:50pm'
$ws.Range("B23").Value = 'Myer, Southland Shopping Centre 1239 Nepean Hwy'
$ws.Range("C23").Value = '22/12/20 10.30am - 11am'
$ws.Range("D23").Value = 'Case visited venue'
$ws.Range("B24").Value = 'Specsavers, 1004-1005 Westfield Southland'
$ws.Range("C24").Value = '22/12/20 11:00am-1145am'
$ws.Range("D24").Value = 'Case shopped in store'
$ws.Range("A25").Value = 'Chirnside Park'
$ws.Range("B25").Value = 'Coles  239-241 Maroondah Hwy, Chirnside Park'
$ws.Range("C25").Value = '31/12/2020 10:00am - 10:15am'
$ws.Range("D25").Value = 'Case shopped at venue'
$ws.Range("A26").Value = 'Clayton'
$ws.Range("B26").Value = 'Kmart - 2107 Dandenong Road, Clayton'
$ws.Range("C26").Value = '30/12/20 7pm - 7.30pm'
$ws.Range("D26").Value = 'Case shopped at store'
$ws.Range("B27").Value = 'Woolworths - M-City, 2107 Dandenong Road, Clayton'
$ws.Range("C27").Value = '30/12/20 7.30pm - 745pm'
$ws.Range("A28").Value = 'Emerald'
$ws.Range("B28").Value = 'Lakeside Paddle Boats, Emerald Lake Park'
$ws.Range("C28").Value = '31/12/20 3:30pm - 5:30pm'
$ws.Range("D28").Value = 'Case visited venue'
$ws.Range("A29").Value = 'Forest Hill'
$ws.Range("B29").Value = 'Forest Hill Chase Shopping Centre 270 Canterbury Rd, Forest Hill VIC 3131'
$ws.Range("C29").Value = '28/12/20 12:00pm-2:00pm'
$ws.Range("D29").Value = '1210hrs Food court 30min; 1250hrs TKMaxx 15min; 1310hrs Target 20min; 1340hrs Woolworths 15min'
$ws.Range("A30").Value = 'Fountain Gate Shopping Centre'
$ws.Range("B30").Value = 'Kmart, Big W, Target, Millers, King of Gifts, Lo Costa  25-55 Overland Drive, Narre Warren VIC 3805'
$ws.Range("C30").Value = '26/12/20 9:00am-11:00am'
$ws.Range("D30").Value = ""
$ws.Range("A31").Value = 'Frankston'
$ws.Range("B31").Value = 'Ishka, Shop G18b, 28 Beach St  Bayside Shopping Centre'
$ws.Range("C31").Value = '31/12/20  3.00pm-3.15pm'
$ws.Range("D31").Value = 'Case shopped'
$ws.Range("B35").Value = 'Sikh Temple Keysborough  198-206 Perry Road, Keysborough'
$ws.Range("C35").Value = '1/01/21 3:00pm-5:00pm'
